# Work on report generation
#
# Removes the placeholder "w" / "Please disregard this..." row from the
# factors sheet, tidies the short_hand value for the Irresponsibility
# Indicator row (drops a stray trailing space: "ii " -> "ii"), and restores
# the last-used cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole placeholder row (was row 18: A="w", C="Please disregard
# this. This does not play any part in the analysis"). Rows below shift up.
$ws.Rows.Item(18).Delete()

# The Irresponsibility Indicator row is now row 18; clean up its short_hand
# value which previously had a trailing space ("ii ").
$ws.Cells.Item(18, 1).Value = "ii"

# Restore the workbook's last active cell selection.
$ws.Range("H11").Select()
